# Generate Report for Handback
# Updates the timestamp values recorded on the handback-status workbook to
# reflect the latest report generation run.

$wb = $excel.ActiveWorkbook

# "Overview" sheet: "Latest HO Xliff Generate Date" column (G) for the
# c8dcfb19-... row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-11-08 23:24:21"

# "zh-cn" sheet: "Correspond Handoff Datetime" (H) and
# "Correspond Handback DateTime" (K) columns for the c8dcfb19-... row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-11-08 23:24:08"
$wsZhCn.Range("K2").Value = "2016-11-08 23:24:44"

# "de-de" sheet: "Correspond Handoff Datetime" (H, shared with the
# "Overview" sheet's "Latest HO Xliff Generate Date") and
# "Correspond Handback DateTime" (K) columns for the c8dcfb19-... row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-11-08 23:24:21"
$wsDeDe.Range("K2").Value = "2016-11-08 23:25:01"
